$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Australia ALeague")

# Rows 159 and 160 had their contents (columns B through AB) swapped,
# while column A (the running index) stayed put.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

foreach ($col in $cols) {
    $addr159 = "${col}159"
    $addr160 = "${col}160"
    $v159 = $ws.Range($addr159).Value2
    $v160 = $ws.Range($addr160).Value2
    $ws.Range($addr159).Value2 = $v160
    $ws.Range($addr160).Value2 = $v159
}

# Row 164 isolated odds updates
$ws.Range("N164").Value = 5
$ws.Range("Q164").Value = 1.89
$ws.Range("R164").Value = 2.01
$ws.Range("T164").Value = 1.875
$ws.Range("U164").Value = 1.975

# Row 165 isolated odds updates
$ws.Range("M165").Value = 2.45
$ws.Range("O165").Value = 2.7
$ws.Range("Q165").Value = 1.85
$ws.Range("R165").Value = 2.05
$ws.Range("T165").Value = 1.8
$ws.Range("U165").Value = 2.05
